$wb = $excel.ActiveWorkbook

# Sheet: 展览 (index 1)
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value = 3304
$ws.Cells.Item(5, 6).Value = 1352
$ws.Cells.Item(6, 6).Value = 43
$ws.Cells.Item(7, 6).Value = 387
$ws.Cells.Item(8, 6).Value = 188
$ws.Cells.Item(10, 6).Value = 8406
$ws.Cells.Item(11, 6).Value = 455
$ws.Cells.Item(13, 6).Value = 82
$ws.Cells.Item(14, 6).Value = 279
$ws.Cells.Item(15, 6).Value = 308
$ws.Cells.Item(16, 6).Value = 127
$ws.Cells.Item(17, 6).Value = 6
$ws.Cells.Item(18, 6).Value = 342
$ws.Cells.Item(19, 6).Value = 10594
$ws.Cells.Item(20, 6).Value = 31
$ws.Cells.Item(23, 6).Value = 25
$ws.Cells.Item(24, 6).Value = 37
$ws.Cells.Item(25, 6).Value = 139
$ws.Cells.Item(27, 6).Value = 180
$ws.Cells.Item(28, 6).Value = 157
$ws.Cells.Item(29, 6).Value = 86
$ws.Cells.Item(31, 6).Value = 92
$ws.Cells.Item(32, 6).Value = 2072
$ws.Cells.Item(33, 6).Value = 34
$ws.Cells.Item(34, 6).Value = 37
$ws.Cells.Item(35, 6).Value = 879
$ws.Cells.Item(37, 6).Value = 273
$ws.Cells.Item(38, 6).Value = 2571
$ws.Cells.Item(39, 6).Value = 3013
$ws.Cells.Item(40, 6).Value = 1228
$ws.Cells.Item(42, 6).Value = 751
$ws.Cells.Item(43, 6).Value = 65
$ws.Cells.Item(44, 6).Value = 327
$ws.Cells.Item(45, 6).Value = 281
$ws.Cells.Item(46, 6).Value = 34
$ws.Cells.Item(47, 6).Value = 100
$ws.Cells.Item(48, 6).Value = 83
$ws.Cells.Item(49, 6).Value = 86
$ws.Cells.Item(50, 6).Value = 67

# Sheet: 演出 (index 2)
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 6).Value = 12
$ws.Cells.Item(9, 6).Value = 50
$ws.Cells.Item(14, 6).Value = 34
$ws.Cells.Item(15, 6).Value = 8
$ws.Cells.Item(16, 6).Value = 40
$ws.Cells.Item(18, 6).Value = 173
$ws.Cells.Item(22, 6).Value = 40
$ws.Cells.Item(23, 6).Value = 23

# Sheet: 本地生活 (index 3)
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(3, 6).Value = 11

# Sheet: 全部类型 (index 4)
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 3304
$ws.Cells.Item(3, 6).Value = 186
$ws.Cells.Item(5, 6).Value = 1352
$ws.Cells.Item(6, 6).Value = 387
$ws.Cells.Item(8, 6).Value = 50
$ws.Cells.Item(9, 6).Value = 188
$ws.Cells.Item(10, 6).Value = 50
$ws.Cells.Item(11, 6).Value = 8406
$ws.Cells.Item(12, 6).Value = 455
$ws.Cells.Item(14, 6).Value = 78
$ws.Cells.Item(15, 6).Value = 82
$ws.Cells.Item(16, 6).Value = 279
$ws.Cells.Item(17, 6).Value = 308
$ws.Cells.Item(18, 6).Value = 6
$ws.Cells.Item(19, 6).Value = 342
$ws.Cells.Item(20, 6).Value = 10594
$ws.Cells.Item(21, 6).Value = 31
$ws.Cells.Item(22, 6).Value = 281
$ws.Cells.Item(25, 6).Value = 139
$ws.Cells.Item(26, 6).Value = 387
$ws.Cells.Item(27, 6).Value = 180
$ws.Cells.Item(28, 6).Value = 34
$ws.Cells.Item(29, 6).Value = 157
$ws.Cells.Item(30, 6).Value = 86
$ws.Cells.Item(31, 6).Value = 2072
$ws.Cells.Item(33, 6).Value = 37
$ws.Cells.Item(34, 6).Value = 879
$ws.Cells.Item(35, 6).Value = 173
$ws.Cells.Item(36, 6).Value = 273
$ws.Cells.Item(37, 6).Value = 2571
$ws.Cells.Item(38, 6).Value = 3013
$ws.Cells.Item(39, 6).Value = 1228
$ws.Cells.Item(41, 6).Value = 751
$ws.Cells.Item(42, 6).Value = 65
$ws.Cells.Item(43, 6).Value = 327
$ws.Cells.Item(44, 6).Value = 23
$ws.Cells.Item(45, 6).Value = 281
$ws.Cells.Item(46, 6).Value = 34
$ws.Cells.Item(47, 6).Value = 100
$ws.Cells.Item(48, 6).Value = 83
$ws.Cells.Item(49, 6).Value = 86
$ws.Cells.Item(50, 6).Value = 67
